$wb = $excel.ActiveWorkbook

# Sheet: ALC (index 1)
$ws = $wb.Worksheets.Item(1)

$ws.Range("H15").Value = 1956.5287
$ws.Range("I15").Value = 1956.5287
$ws.Range("K15").Value = 5869.5861
$ws.Range("M15").Value = -5700.5861
$ws.Range("H17").Value = 388240.06
$ws.Range("I17").Value = 199
$ws.Range("J17").Value = 455725.44
$ws.Range("K17").Value = 597
$ws.Range("L17").Value = 1367176.32
$ws.Range("M17").Value = -429
$ws.Range("N17").Value = -1367512.32
$ws.Range("H33").Value = 133.66667
$ws.Range("I33").Value = 150
$ws.Range("K33").Value = 150
$ws.Range("M33").Value = 79
$ws.Range("H116").Value = 2366.36
$ws.Range("I116").Value = 2474.4546
$ws.Range("J116").Value = 2281.4285
$ws.Range("K116").Value = 2474.4546
$ws.Range("L116").Value = 2281.4285
$ws.Range("M116").Value = 967.5454
$ws.Range("N116").Value = -9165.4285
$ws.Range("H132").Value = 2781.8906
$ws.Range("I132").Value = 2285.3555
$ws.Range("J132").Value = 3957.8948
$ws.Range("K132").Value = 6856.066500000001
$ws.Range("L132").Value = 11873.6844
$ws.Range("M132").Value = -4326.066500000001
$ws.Range("N132").Value = -16933.6844
$ws.Range("H137").Value = 3929.449
$ws.Range("I137").Value = 1639
$ws.Range("J137").Value = 6315.3335
$ws.Range("K137").Value = 4917
$ws.Range("L137").Value = 18946.0005
$ws.Range("M137").Value = -2367
$ws.Range("N137").Value = -24046.0005
$ws.Range("H138").Value = 2301.9153
$ws.Range("I138").Value = 1235.6052
$ws.Range("J138").Value = 4231.4287
$ws.Range("K138").Value = 3706.8156
$ws.Range("L138").Value = 12694.2861
$ws.Range("M138").Value = 1433.1844
$ws.Range("N138").Value = -22974.2861
# Sheet: ARM (index 2)
$ws = $wb.Worksheets.Item(2)

$ws.Range("H32").Value = 11770046
$ws.Range("I32").Value = 12662508
$ws.Range("J32").Value = 19300
$ws.Range("K32").Value = 12662508
$ws.Range("L32").Value = 19300
$ws.Range("M32").Value = -12662221
$ws.Range("N32").Value = -19874
$ws.Range("H61").Value = 1377.1464
$ws.Range("I61").Value = 1139.8438
$ws.Range("K61").Value = 1139.8438
$ws.Range("M61").Value = -927.8438000000001
$ws.Range("H74").Value = 2746.1458
$ws.Range("I74").Value = 507.66666
$ws.Range("J74").Value = 12446.223
$ws.Range("K74").Value = 507.66666
$ws.Range("L74").Value = 12446.223
$ws.Range("M74").Value = 366.33334
$ws.Range("N74").Value = -14194.223
$ws.Range("H77").Value = 2746.1458
$ws.Range("I77").Value = 507.66666
$ws.Range("J77").Value = 12446.223
$ws.Range("K77").Value = 2538.3333
$ws.Range("L77").Value = 62231.115
$ws.Range("M77").Value = 1829.6667
$ws.Range("N77").Value = -70967.11499999999
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").Value = ""
$ws.Range("H132").Value = 8230.645500000001
$ws.Range("I132").Value = 6295.7144
$ws.Range("J132").Value = 12294
$ws.Range("K132").Value = 18887.1432
$ws.Range("L132").Value = 36882
$ws.Range("M132").Value = -16357.1432
$ws.Range("N132").Value = -41942
$ws.Range("H136").Value = 1377.1464
$ws.Range("I136").Value = 1139.8438
$ws.Range("K136").Value = 3419.5314
$ws.Range("M136").Value = -869.5314000000003
# Sheet: BSM (index 3)
$ws = $wb.Worksheets.Item(3)

$ws.Range("H80").Value = 505.94446
$ws.Range("I80").Value = 630.4286
$ws.Range("J80").Value = 426.72726
$ws.Range("K80").Value = 630.4286
$ws.Range("L80").Value = 426.72726
$ws.Range("M80").Value = 367.5714
$ws.Range("N80").Value = -2422.72726
$ws.Range("H83").Value = 505.94446
$ws.Range("I83").Value = 630.4286
$ws.Range("J83").Value = 426.72726
$ws.Range("K83").Value = 3152.143
$ws.Range("L83").Value = 2133.6363
$ws.Range("M83").Value = 1839.857
$ws.Range("N83").Value = -12117.6363
$ws.Range("H86").Value = 1864.6487
$ws.Range("I86").Value = 1601.1428
$ws.Range("J86").Value = 2210.5
$ws.Range("K86").Value = 1601.1428
$ws.Range("L86").Value = 2210.5
$ws.Range("M86").Value = -478.1428000000001
$ws.Range("N86").Value = -4456.5
$ws.Range("H89").Value = 1864.6487
$ws.Range("I89").Value = 1601.1428
$ws.Range("J89").Value = 2210.5
$ws.Range("K89").Value = 8005.714
$ws.Range("L89").Value = 11052.5
$ws.Range("M89").Value = -2389.714
$ws.Range("N89").Value = -22284.5
# Sheet: CRP (index 4)
$ws = $wb.Worksheets.Item(4)

$ws.Range("H31").Value = 19651774
$ws.Range("I31").Value = 38463030
$ws.Range("J31").Value = 88062.8
$ws.Range("K31").Value = 38463030
$ws.Range("L31").Value = 88062.8
$ws.Range("M31").Value = -38462735
$ws.Range("N31").Value = -88652.8
$ws.Range("H34").Value = 19651774
$ws.Range("I34").Value = 38463030
$ws.Range("J34").Value = 88062.8
$ws.Range("K34").Value = 38463030
$ws.Range("L34").Value = 88062.8
$ws.Range("M34").Value = -38462828
$ws.Range("N34").Value = -88466.8
$ws.Range("H58").Value = 969.8701
$ws.Range("I58").Value = 664.96
$ws.Range("J58").Value = 1534.5186
$ws.Range("K58").Value = 664.96
$ws.Range("L58").Value = 1534.5186
$ws.Range("M58").Value = -461.96
$ws.Range("N58").Value = -1940.5186
$ws.Range("H105").Value = 2149.577
$ws.Range("I105").Value = 1414.6471
$ws.Range("J105").Value = 3537.7778
$ws.Range("K105").Value = 1414.6471
$ws.Range("L105").Value = 3537.7778
$ws.Range("M105").Value = 332.3529000000001
$ws.Range("N105").Value = -7031.7778
$ws.Range("H107").Value = 46250.047
$ws.Range("I107").Value = 72212.71000000001
$ws.Range("J107").Value = 815.375
$ws.Range("K107").Value = 72212.71000000001
$ws.Range("L107").Value = 815.375
$ws.Range("M107").Value = -70292.71000000001
$ws.Range("N107").Value = -4655.375
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").Value = ""
$ws.Range("H120").Value = 11071.429
$ws.Range("J120").Value = 11071.429
$ws.Range("L120").Value = 11071.429
$ws.Range("N120").Value = -18329.429
$ws.Range("H121").Value = 24400
$ws.Range("J121").Value = 24400
$ws.Range("L121").Value = 24400
$ws.Range("N121").Value = -27020
$ws.Range("H132").Value = 15628882
$ws.Range("I132").Value = 20411778
$ws.Range("J132").Value = 4752.4
$ws.Range("K132").Value = 61235334
$ws.Range("L132").Value = 14257.2
$ws.Range("M132").Value = -61232804
$ws.Range("N132").Value = -19317.2
$ws.Range("H134").Value = 2015.6216
$ws.Range("I134").Value = 1930.125
$ws.Range("J134").Value = 2562.8
$ws.Range("K134").Value = 5790.375
$ws.Range("L134").Value = 7688.400000000001
$ws.Range("M134").Value = -3255.375
$ws.Range("N134").Value = -12758.4
$ws.Range("H136").Value = 969.8701
$ws.Range("I136").Value = 664.96
$ws.Range("J136").Value = 1534.5186
$ws.Range("K136").Value = 1994.88
$ws.Range("L136").Value = 4603.5558
$ws.Range("M136").Value = 555.1199999999999
$ws.Range("N136").Value = -9703.5558
# Sheet: CUL (index 5)
$ws = $wb.Worksheets.Item(5)

$ws.Range("H129").Value = 771.7692
$ws.Range("I129").Value = 750
$ws.Range("J129").Value = 844.3333
$ws.Range("K129").Value = 2250
$ws.Range("L129").Value = 2532.9999
$ws.Range("M129").Value = 2750
$ws.Range("N129").Value = -12532.9999
# Sheet: GSM (index 6)
$ws = $wb.Worksheets.Item(6)

$ws.Range("H42").Value = 39999.168
$ws.Range("J42").Value = 39999.168
$ws.Range("L42").Value = 39999.168
$ws.Range("N42").Value = -40969.168
$ws.Range("H115").Value = 39999.168
$ws.Range("J115").Value = 39999.168
$ws.Range("L115").Value = 39999.168
$ws.Range("N115").Value = -42349.168
$ws.Range("H117").Value = 29000
$ws.Range("J117").Value = 29000
$ws.Range("L117").Value = 29000
$ws.Range("N117").Value = -35884
$ws.Range("H119").Value = 48000
$ws.Range("J119").Value = 48000
$ws.Range("L119").Value = 48000
$ws.Range("N119").Value = -57676
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").Value = ""
$ws.Range("H132").Value = 3599.228
$ws.Range("I132").Value = 3756.2341
$ws.Range("J132").Value = 2861.3
$ws.Range("K132").Value = 11268.7023
$ws.Range("L132").Value = 8583.900000000001
$ws.Range("M132").Value = -8738.702300000001
$ws.Range("N132").Value = -13643.9
# Sheet: LTW (index 7)
$ws = $wb.Worksheets.Item(7)

$ws.Range("H115").Value = 41900
$ws.Range("J115").Value = 41900
$ws.Range("L115").Value = 41900
$ws.Range("N115").Value = -44250
$ws.Range("H117").Value = 41000
$ws.Range("J117").Value = 41000
$ws.Range("L117").Value = 41000
$ws.Range("N117").Value = -50178
$ws.Range("H119").Value = 33557.145
$ws.Range("J119").Value = 33557.145
$ws.Range("L119").Value = 33557.145
$ws.Range("N119").Value = -43233.145
$ws.Range("H120").Value = 42750
$ws.Range("J120").Value = 42750
$ws.Range("L120").Value = 42750
$ws.Range("N120").Value = -52426
$ws.Range("H136").Value = 3098.018
$ws.Range("I136").Value = 1218.65
$ws.Range("J136").Value = 8109.6665
$ws.Range("K136").Value = 3655.95
$ws.Range("L136").Value = 24328.9995
$ws.Range("M136").Value = -1105.95
$ws.Range("N136").Value = -29428.9995
# Sheet: WVR (index 8)
$ws = $wb.Worksheets.Item(8)

$ws.Range("H81").Value = 967.55
$ws.Range("I81").Value = 926.5294
$ws.Range("J81").Value = 1200
$ws.Range("K81").Value = 1853.0588
$ws.Range("L81").Value = 2400
$ws.Range("M81").Value = -792.0588
$ws.Range("N81").Value = -4522
$ws.Range("H84").Value = 967.55
$ws.Range("I84").Value = 926.5294
$ws.Range("J84").Value = 1200
$ws.Range("K84").Value = 9265.294
$ws.Range("L84").Value = 12000
$ws.Range("M84").Value = -3961.294
$ws.Range("N84").Value = -22608
$ws.Range("H116").Value = 48000
$ws.Range("J116").Value = 48000
$ws.Range("L116").Value = 48000
$ws.Range("N116").Value = -57178
$ws.Range("H117").Value = 37803
$ws.Range("J117").Value = 37803
$ws.Range("L117").Value = 37803
$ws.Range("N117").Value = -46981
$ws.Range("H132").Value = 8199531.5
$ws.Range("I132").Value = 12824004
$ws.Range("J132").Value = 1602.4546
$ws.Range("K132").Value = 38472012
$ws.Range("L132").Value = 4807.3638
$ws.Range("M132").Value = -38469482
$ws.Range("N132").Value = -9867.363799999999
$ws.Range("H136").Value = 1036.519
$ws.Range("I136").Value = 504.16666
$ws.Range("J136").Value = 1640.8108
$ws.Range("K136").Value = 1512.49998
$ws.Range("L136").Value = 4922.4324
$ws.Range("M136").Value = 1037.50002
$ws.Range("N136").Value = -10022.4324
